$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MS2_charges")

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "1,2,3,4,5,6,7,8"

$ws.Range("B21").Select()
